$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orgs = @(
    "Owl Ventures",
    "Bill & Melinda Gates Foundation",
    "Penn State University - Outreach & Engagement",
    "Chan Zuckerberg Initiative (CZI)",
    "Google for Education",
    "Y Combinator",
    "National Science Foundation (NSF)",
    "Reach Capital",
    "U.S. Dept. of Education (EIR Program)",
    "IES SBIR (ED/IES)",
    "500 Global Flagship VC (non-accelerator checks)",
    "TGR Foundation (Tiger Woods)",
    "Duke University - Duke-Durham Neighborhood Partnership",
    "500 Global (seed/accelerator)",
    "The Ohio State University - Office of Outreach",
    "New York Knicks Garden of Dreams Foundation",
    "Miami Dolphins Foundation",
    "Los Angeles Lakers Youth Foundation",
    "Seattle Seahawks - Spirit of 12",
    "Kansas City Chiefs Foundation",
    "Chicago Bears Charities",
    "Chicago Cubs Charities",
    "Los Angeles Rams Foundation",
    "Toronto Maple Leafs - MLSE Foundation",
    "FC Dallas Foundation",
    "Minnesota Wild Foundation",
    "Houston Texans Foundation",
    "Jacksonville Jaguars Foundation",
    "Tampa Bay Lightning Community Heroes",
    "Laureus Sport for Good USA"
)

for ($i = 0; $i -lt $orgs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $orgs[$i]
}
